$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.527.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.911.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4845'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.08%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06805'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '111.07'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.29'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.908.68'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07561'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.400'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '294.96'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.524.33'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.03'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007585'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.511'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.160.33'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.423'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.461'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.32'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.079'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1063'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.137'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04974'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7360'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.681'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.021'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.16'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4443'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8653'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.766'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9993'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '69.24'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.206'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.24'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.189'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1228'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2515'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.01%  '
